$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New accelerometer dataset (x, y, z) replacing the old rows 2-22 with new rows 2-21.
# The first 7 rows correspond to the old rows 16-22 (the struggle class data that
# was already present), followed by 13 brand-new rows for the newly added
# falling / walkingToRunning classes, per the commit message.
$data = @(
    @(-0.2141320705413827, -2.326400130987168, -2.468529105186464),
    @(2.301481246948242, -1.784020185470581, 0.5732603073120117),
    @(0.7297354936599727, -2.234472751617432, -4.658630132675173),
    @(1.63190019130707, -2.561783850193024, -1.591027021408075),
    @(0.6863539814949051, -1.503557689487935, 1.746999144554136),
    @(0.2367095947265618, -1.311740666627883, -0.07715380191802312),
    @(1.296695142984393, -3.332512527704244, -2.418631196022031),
    @(1.87941366434097, -4.668229699134821, -1.045372545719146),
    @(0.3820920586585984, -1.311929136514663, 0.04638075828551969),
    @(2.75743055343628, -3.144901037216187, 4.137303829193115),
    @(4.460695505142212, 1.830066174268722, -0.1949661374092102),
    @(-3.40113162994386, 3.317261695861818, 1.794482350349431),
    @(0.008035421371476836, -0.4300747811794388, 2.762799173593522),
    @(1.565377473831186, -1.947239398956303, 1.930309116840367),
    @(-0.791193664073943, -1.430967807769773, 0.315328881144527),
    @(0.01700598001480103, -0.8887928128242493, 0.1178494691848755),
    @(3.75538071990013, 1.982138156890872, 1.528477013111122),
    @(4.293412685394287, -2.266220092773437, -9.760974884033203),
    @(2.267539381980895, -2.359471559524536, -1.19436234235763),
    @(-1.645043730735778, -1.020936071872711, 0.8753915429115292)
)

$oldLastRow = 22
$newLastRow = 1 + $data.Count

# Clear out the old data range entirely first (including the row that will
# disappear from the used range), then write the new values.
$ws.Range("A2:C$oldLastRow").ClearContents()

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
